$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table refresh (scheduled GitHub Actions run).
# Column D ("Price") values are stored as plain text even when they look like
# numbers (e.g. locale-formatted "67.070.28" using "." as a thousands separator,
# or plain decimals like "1.01"). Force the whole Price column to text format
# first so Excel does not silently reinterpret the new values as numbers, then
# restore the original (default) style once all values are written.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range('D2').Value = '67.178.71'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '2.614.94'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').Value = '593.49'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = '152.55'
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').Value = '2.611.32'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('D12').Value = '5.19'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('E13').Value = '  -3.62%  '
$ws.Range('D14').Value = '27.62'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').Value = '3.087.24'
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('E16').Value = '  -4.51%  '
$ws.Range('D17').Value = '67.019.40'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').Value = '2.606.15'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = '364.96'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').Value = '11.04'
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('D21').Value = '7.34'
$ws.Range('E21').Value = '  -5.25%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').Value = '2.06'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D25').Value = '10.10'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').Value = '67.49'
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('D27').Value = '2.741.49'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').Value = '587.34'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  -4.16%  '
$ws.Range('E31').Value = '  -5.06%  '
$ws.Range('D32').Value = '7.71'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('E33').Value = '  -2.68%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '0.124'
$ws.Range('E35').Value = '  -6.90%  '
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('D38').Value = '156.21'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '19.03'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = '0.366'
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('D41').Value = '5.24'
$ws.Range('E41').Value = '  -3.44%  '
$ws.Range('D42').Value = '1.81'
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').Value = '40.79'
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '16.49'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '155.03'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').Value = '0.0₆0294'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '3.72'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '21.60'
$ws.Range('E50').Value = '  +2.62%  '
$ws.Range('D51').Value = '0.617'
$ws.Range('E51').Value = '  -3.24%  '

$priceCol.Style = "Normal"
